$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = 37.31870588235294

for ($row = 24; $row -le 42; $row++) {
    $ws.Cells.Item($row, 9).Value = $newValue
}
